# Fruta / hortaliza, semanal
# New weekly price observations for "Plátano" (Pintón / Primera Pintón)
# at "Terminal Hortofrutícola Agro Chillán" are inserted at the top of the
# date-ordered (descending) block, pushing the existing rows 773-804 down
# to 775-806 (dimension grows from A1:T804 to A1:T806).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 773-774; everything below shifts down by 2.
$ws.Rows("773:774").Insert()

# New row 773: Plátano, Pintón
$ws.Range("A773").Value = 7
$ws.Range("B773").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C773").Value = "Ñuble"
$ws.Range("D773").Value = 44939
$ws.Range("E773").Value = 16
$ws.Range("F773").Value = "Fruta"
$ws.Range("G773").Value = 100108
$ws.Range("H773").Value = "Tropicales y subtropicales"
$ws.Range("I773").Value = 100108006
$ws.Range("J773").Value = "Plátano"
$ws.Range("K773").Value = "Sin especificar"
$ws.Range("L773").Value = "Pintón"
$ws.Range("M773").Value = 250
$ws.Range("N773").Value = 22000
$ws.Range("O773").Value = 22000
$ws.Range("P773").Value = 22000
$ws.Range("Q773").Value = "$/caja 20 kilos"
$ws.Range("R773").Value = "Ecuador"
$ws.Range("S773").Value = 1100
$ws.Range("T773").Value = 20

# New row 774: Plátano, Primera Pintón
$ws.Range("A774").Value = 7
$ws.Range("B774").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C774").Value = "Ñuble"
$ws.Range("D774").Value = 44939
$ws.Range("E774").Value = 16
$ws.Range("F774").Value = "Fruta"
$ws.Range("G774").Value = 100108
$ws.Range("H774").Value = "Tropicales y subtropicales"
$ws.Range("I774").Value = 100108006
$ws.Range("J774").Value = "Plátano"
$ws.Range("K774").Value = "Sin especificar"
$ws.Range("L774").Value = "Primera Pintón"
$ws.Range("M774").Value = 300
$ws.Range("N774").Value = 23000
$ws.Range("O774").Value = 24000
$ws.Range("P774").Value = 23500
$ws.Range("Q774").Value = "$/caja 20 kilos"
$ws.Range("R774").Value = "Ecuador"
$ws.Range("S774").Value = 1175
$ws.Range("T774").Value = 20
